$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add P1 and Q1 headers with same style as O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Row 2
$ws.Range("D2").Value = 2.104931546227315
$ws.Range("E2").Value = 32.08426124340383
$ws.Range("F2").Value = 32.50263912979236
$ws.Range("G2").Value = 2.006822032921141
$ws.Range("H2").Value = 4.685067244831645
$ws.Range("L2").Value = 0
$ws.Range("O2").ClearContents()
$ws.Range("Q2").Value = 25.77189398948176

# Row 3
$ws.Range("D3").Value = 2.07901484599523
$ws.Range("E3").Value = 30.08996734270405
$ws.Range("F3").Value = 30.40737855835052
$ws.Range("G3").Value = 2.014456688235499
$ws.Range("H3").Value = 4.356391808371115
$ws.Range("L3").Value = 0
$ws.Range("O3").ClearContents()
$ws.Range("Q3").Value = 24.17743910735283

# Row 4
$ws.Range("D4").Value = 2.062015653248795
$ws.Range("E4").Value = 28.80488661861066
$ws.Range("F4").Value = 29.06442210359348
$ws.Range("G4").Value = 2.019259790255377
$ws.Range("H4").Value = 4.148825185931026
$ws.Range("L4").Value = 0
$ws.Range("O4").ClearContents()
$ws.Range("Q4").Value = 23.1583291308764

# Row 5
$ws.Range("D5").Value = 2.054817360627801
$ws.Range("E5").Value = 28.2659412095414
$ws.Range("F5").Value = 28.50310277912704
$ws.Range("G5").Value = 2.021247547256435
$ws.Range("H5").Value = 4.062737089695345
$ws.Range("L5").Value = 0
$ws.Range("O5").ClearContents()
$ws.Range("Q5").Value = 22.73311745795231

# Row 6
$ws.Range("D6").Value = 2.053605789454675
$ws.Range("E6").Value = 28.17552759242295
$ws.Range("F6").Value = 28.4090531159221
$ws.Range("G6").Value = 2.021579492430742
$ws.Range("H6").Value = 4.04834990266911
$ws.Range("L6").Value = 0
$ws.Range("O6").ClearContents()
$ws.Range("Q6").Value = 22.66191882739718

# Row 7
$ws.Range("D7").Value = 2.061919668784863
$ws.Range("E7").Value = 28.79768003346285
$ws.Range("F7").Value = 29.05690860914343
$ws.Range("G7").Value = 2.019286472699942
$ws.Range("H7").Value = 4.147670327775133
$ws.Range("L7").Value = 0
$ws.Range("O7").ClearContents()
$ws.Range("Q7").Value = 23.152634451672

# Row 8
$ws.Range("D8").Value = 2.096222154273882
$ws.Range("E8").Value = 31.4091702768249
$ws.Range("F8").Value = 31.79197252358496
$ws.Range("G8").Value = 2.009431407985452
$ws.Range("H8").Value = 4.572900976867171
$ws.Range("L8").Value = 0
$ws.Range("O8").ClearContents()
$ws.Range("Q8").Value = 25.23053070064459

# Row 9
$ws.Range("D9").Value = 2.154785196017936
$ws.Range("E9").Value = 36.05303923638011
$ws.Range("F9").Value = 36.83731268649606
$ws.Range("G9").Value = 1.990950926588451
$ws.Range("H9").Value = 5.363758942983138
$ws.Range("L9").Value = 0
$ws.Range("O9").ClearContents()
$ws.Range("Q9").Value = 29.08145777737296

# Row 10
$ws.Range("D10").Value = 2.167284639790014
$ws.Range("E10").Value = 38.30496531119682
$ws.Range("F10").Value = 40.11688766140962
$ws.Range("G10").Value = 1.978295844802251
$ws.Range("H10").Value = 5.874786799873153
$ws.Range("L10").Value = 0
$ws.Range("O10").ClearContents()
$ws.Range("Q10").Value = 31.67190502327665

# Row 11
$ws.Range("D11").Value = 1.988040873747462
$ws.Range("E11").Value = 31.90383516709814
$ws.Range("F11").Value = 39.69607580327018
$ws.Range("G11").Value = 1.976787567993235
$ws.Range("H11").Value = 6.076936552026438
$ws.Range("L11").Value = 0
$ws.Range("O11").ClearContents()
$ws.Range("Q11").Value = 31.16269187899136

# Row 12
$ws.Range("D12").Value = 1.886694621792792
$ws.Range("E12").Value = 26.22082948752583
$ws.Range("F12").Value = 38.66876359871462
$ws.Range("G12").Value = 1.977779670963913
$ws.Range("H12").Value = 6.631181598599258
$ws.Range("L12").Value = 0
$ws.Range("O12").ClearContents()
$ws.Range("Q12").Value = 30.20426100823344

# Row 13
$ws.Range("D13").Value = 1.845175175909834
$ws.Range("E13").Value = 20.74403178807173
$ws.Range("F13").Value = 37.08093892212274
$ws.Range("G13").Value = 1.980732405307392
$ws.Range("H13").Value = 7.402518963799866
$ws.Range("L13").Value = 0
$ws.Range("O13").ClearContents()
$ws.Range("Q13").Value = 28.81279978420219

# Row 14
$ws.Range("D14").Value = 1.85182681988758
$ws.Range("E14").Value = 17.18237530305624
$ws.Range("F14").Value = 35.69629785485314
$ws.Range("G14").Value = 1.98362620208446
$ws.Range("H14").Value = 8.048460287447998
$ws.Range("L14").Value = 0
$ws.Range("O14").ClearContents()
$ws.Range("Q14").Value = 27.62586430547521

# Row 15
$ws.Range("D15").Value = 1.858155568471455
$ws.Range("E15").Value = 16.35316398773954
$ws.Range("F15").Value = 35.23247085823836
$ws.Range("G15").Value = 1.98479689414014
$ws.Range("H15").Value = 8.195877317525778
$ws.Range("L15").Value = 0
$ws.Range("O15").ClearContents()
$ws.Range("Q15").Value = 27.2380093710241

# Row 16
$ws.Range("D16").Value = 1.859260628164319
$ws.Range("E16").Value = 16.10521882252414
$ws.Range("F16").Value = 34.07784629588687
$ws.Range("G16").Value = 1.989491393222659
$ws.Range("H16").Value = 7.892855105775236
$ws.Range("L16").Value = 0
$ws.Range("O16").ClearContents()
$ws.Range("Q16").Value = 26.34182787474889

# Row 17
$ws.Range("D17").Value = 1.851924556190423
$ws.Range("E17").Value = 17.77691526196038
$ws.Range("F17").Value = 33.95121982846675
$ws.Range("G17").Value = 1.991606321526372
$ws.Range("H17").Value = 7.265542507698992
$ws.Range("L17").Value = 0
$ws.Range("O17").ClearContents()
$ws.Range("Q17").Value = 26.30645079517831

# Row 18
$ws.Range("D18").Value = 1.867210366100637
$ws.Range("E18").Value = 21.78599612047764
$ws.Range("F18").Value = 34.71660349521666
$ws.Range("G18").Value = 1.991542102143376
$ws.Range("H18").Value = 6.406190035247825
$ws.Range("L18").Value = 0
$ws.Range("O18").ClearContents()
$ws.Range("Q18").Value = 27.02733415734664

# Row 19
$ws.Range("D19").Value = 1.943077678630588
$ws.Range("E19").Value = 27.6789941233401
$ws.Range("F19").Value = 36.1383136191745
$ws.Range("G19").Value = 1.989428526737569
$ws.Range("H19").Value = 5.706895936690269
$ws.Range("L19").Value = 0
$ws.Range("O19").ClearContents()
$ws.Range("Q19").Value = 28.3002702775884

# Row 20
$ws.Range("D20").Value = 2.16322082627052
$ws.Range("E20").Value = 37.68523828858719
$ws.Range("F20").Value = 39.26505573281266
$ws.Range("G20").Value = 1.981616666940774
$ws.Range("H20").Value = 5.739149298123078
$ws.Range("L20").Value = 0
$ws.Range("O20").ClearContents()
$ws.Range("Q20").Value = 30.99836418970007

# Row 21
$ws.Range("D21").Value = 2.210031903568887
$ws.Range("E21").Value = 40.69173931999269
$ws.Range("F21").Value = 42.02310993126699
$ws.Range("G21").Value = 1.971206310666587
$ws.Range("H21").Value = 6.200291860830716
$ws.Range("L21").Value = 0
$ws.Range("O21").ClearContents()
$ws.Range("Q21").Value = 33.20735580540065

# Row 22
$ws.Range("D22").Value = 2.226764471669936
$ws.Range("E22").Value = 42.15519950247793
$ws.Range("F22").Value = 43.66009698774357
$ws.Range("G22").Value = 1.964699594393956
$ws.Range("H22").Value = 6.474213141561632
$ws.Range("L22").Value = 0
$ws.Range("O22").ClearContents()
$ws.Range("Q22").Value = 34.50849497634176

# Row 23
$ws.Range("D23").Value = 2.217925893142396
$ws.Range("E23").Value = 41.3784766385451
$ws.Range("F23").Value = 42.79119724741173
$ws.Range("G23").Value = 1.968169896246649
$ws.Range("H23").Value = 6.328223990864769
$ws.Range("L23").Value = 0
$ws.Range("O23").ClearContents()
$ws.Range("Q23").Value = 33.81792256225349

# Row 24
$ws.Range("D24").Value = 2.182435951789251
$ws.Range("E24").Value = 38.3397785689396
$ws.Range("F24").Value = 39.39336908754179
$ws.Range("G24").Value = 1.981382213315997
$ws.Range("H24").Value = 5.770189266003064
$ws.Range("L24").Value = 0
$ws.Range("O24").ClearContents()
$ws.Range("Q24").Value = 31.11596966726311

# Row 25
$ws.Range("D25").Value = 2.139885247190772
$ws.Range("E25").Value = 34.84722144548751
$ws.Range("F25").Value = 35.48956984981111
$ws.Range("G25").Value = 1.995870373675259
$ws.Range("H25").Value = 5.154018993962229
$ws.Range("L25").Value = 0
$ws.Range("O25").ClearContents()
$ws.Range("Q25").Value = 28.00786865070119
